$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking price cells as Text so Excel keeps them as
# literal strings (matching the source data) instead of coercing to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "43.079.92"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.397.16"
$ws.Range("E3").Value = "  +5.10%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "333.65"
$ws.Range("E5").Value = "  +9.55%  "
$ws.Range("D6").Value = "105.69"
$ws.Range("E6").Value = "  -7.33%  "
$ws.Range("D7").Value = "0.651"
$ws.Range("E7").Value = "  +3.11%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "0.650"
$ws.Range("E9").Value = "  +5.72%  "
$ws.Range("D10").Value = "42.19"
$ws.Range("E10").Value = "  -5.71%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "8.73"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "17.11"
$ws.Range("E14").Value = "  +11.04%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "2.756.88"
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("D17").Value = "2.398.74"
$ws.Range("E17").Value = "  +5.09%  "
$ws.Range("D18").Value = "43.111.35"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +6.54%  "
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "3.83"
$ws.Range("E21").Value = "  +7.98%  "
$ws.Range("D22").Value = "77.33"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "275.15"
$ws.Range("E23").Value = "  +8.00%  "
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").Value = "9.80"
$ws.Range("E25").Value = "  +8.59%  "
$ws.Range("D26").Value = "11.94"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "23.26"
$ws.Range("E28").Value = "  +4.78%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "174.82"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "37.16"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "0.0937"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("D34").Value = "6.02"
$ws.Range("E34").Value = "  +5.53%  "
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").Value = "4.90"
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("D37").Value = "4.11"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("D38").Value = "0.0365"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +11.26%  "
$ws.Range("E41").Value = "  +13.41%  "
$ws.Range("D42").Value = "0.236"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").Value = "70.16"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").Value = "122.14"
$ws.Range("E44").Value = "  +14.41%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "91.93"
$ws.Range("E46").Value = "  +44.71%  "
$ws.Range("D47").Value = "12.35"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").Value = "5.58"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  +5.75%  "
$ws.Range("D50").Value = "0.514"
$ws.Range("E50").Value = "  +17.72%  "
$ws.Range("E51").Value = "  +1.98%  "
